$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Förändrad" (Changed) date column C, for all data rows (2 through 28),
# is updated from serial date 45462 (2024-06-19) to 45464 (2024-06-21).
for ($row = 2; $row -le 28; $row++) {
    $ws.Range("C$row").Value = 45464
}
